# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "Datos actualizados" timestamp title (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 03:35"

# --- Bolivia overtakes Finlandia in the ranking (rows 66 & 67) ---
# Row 66 keeps pointing at the country that is now ranked #70 (Bolivia),
# and gets freshly updated statistics.
$ws.Range("A66").Value = "Bolivia"
$ws.Range("B66").Value = 6660
$ws.Range("C66").Value = 397
$ws.Range("D66").Value = 647
$ws.Range("E66").Value = 5752
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 11
$ws.Range("H66").Value = 261

# Row 67 now points at the country ranked #71 (Finlandia), which carries
# the statistics that used to belong to row 66 (Finlandia's old figures).
$ws.Range("A67").Value = "Finlandia"
$ws.Range("B67").Value = 6599
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 5100
$ws.Range("E67").Value = 1191
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 308

# --- Independent data refreshes ---

# Row 49
$ws.Range("B49").Value = 11225
$ws.Range("C49").Value = 19
$ws.Range("D49").Value = 10275
$ws.Range("E49").Value = 681
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 269

# Row 73
$ws.Range("B73").Value = 3976
$ws.Range("C73").Value = 156
$ws.Range("D73").Value = 503
$ws.Range("E73").Value = 3303
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 170

# Row 98 (only active/recovered split changes)
$ws.Range("D98").Value = 1461
$ws.Range("E98").Value = 22
